$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells retain text formatting so purely-numeric
# looking strings (e.g. "1.005") are not reinterpreted as numbers by Excel.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

$ws.Range('D2').Value = '27.238.26'
$ws.Range('E2').Value = '  -1.90%  '
$ws.Range('D3').Value = '1.821.80'
$ws.Range('E3').Value = '  -1.93%  '
$ws.Range('D4').Value = '1.005'
$ws.Range('E4').Value = '  -1.30%  '
$ws.Range('D5').Value = '314.16'
$ws.Range('E5').Value = '  -2.01%  '
$ws.Range('E6').Value = '  -1.38%  '
$ws.Range('D7').Value = '0.4278'
$ws.Range('E7').Value = '  -2.15%  '
$ws.Range('D8').Value = '0.3680'
$ws.Range('E8').Value = '  -2.59%  '
$ws.Range('B9').Value = 'Dogecoin'
$ws.Range('C9').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D9').Value = '0.07234'
$ws.Range('E9').Value = '  -2.63%  '
$ws.Range('B10').Value = 'Polygon'
$ws.Range('C10').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D10').Value = '0.8617'
$ws.Range('E10').Value = '  -2.49%  '
$ws.Range('B11').Value = 'Solana'
$ws.Range('C11').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D11').Value = '20.99'
$ws.Range('E11').Value = '  -2.59%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.829.49'
$ws.Range('E12').Value = '  -1.62%  '
$ws.Range('B13').Value = 'Chainlink'
$ws.Range('C13').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D13').Value = '6.673'
$ws.Range('E13').Value = '  -1.15%  '
$ws.Range('B14').Value = 'TRON'
$ws.Range('C14').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D14').Value = '0.07134'
$ws.Range('E14').Value = '  +0.47%  '
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').Value = '5.306'
$ws.Range('E15').Value = '  -3.38%  '
$ws.Range('B16').Value = 'Litecoin'
$ws.Range('C16').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D16').Value = '89.20'
$ws.Range('E16').Value = '  +1.19%  '
$ws.Range('B17').Value = 'BinanceUSD'
$ws.Range('C17').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D17').Value = '1.006'
$ws.Range('E17').Value = '  -1.61%  '
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').Value = '0.000008862'
$ws.Range('E18').Value = '  -1.84%  '
$ws.Range('B19').Value = 'Dai'
$ws.Range('C19').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D19').Value = '1.004'
$ws.Range('E19').Value = '  -1.33%  '
$ws.Range('B20').Value = 'Avalanche'
$ws.Range('C20').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D20').Value = '15.05'
$ws.Range('E20').Value = '  -2.72%  '
$ws.Range('B21').Value = 'WrappedBTC'
$ws.Range('C21').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D21').Value = '27.258.76'
$ws.Range('E21').Value = '  -1.82%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').Value = '5.148'
$ws.Range('E22').Value = '  -2.44%  '
$ws.Range('B23').Value = 'Cosmos'
$ws.Range('C23').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D23').Value = '10.89'
$ws.Range('E23').Value = '  -2.39%  '
$ws.Range('B24').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C24').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D24').Value = '2.051.92'
$ws.Range('E24').Value = '  -1.77%  '
$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').Value = '2.006'
$ws.Range('E25').Value = '  -1.47%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').Value = '153.29'
$ws.Range('E26').Value = '  -2.35%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = '18.33'
$ws.Range('E27').Value = '  -1.86%  '
$ws.Range('B28').Value = 'LidoDAOToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D28').Value = '2.126'
$ws.Range('E28').Value = '  +6.67%  '
$ws.Range('B29').Value = 'InternetComputer(DFINITY)'
$ws.Range('C29').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D29').Value = '5.233'
$ws.Range('E29').Value = '  -3.89%  '
$ws.Range('B30').Value = 'BitcoinCash'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D30').Value = '116.47'
$ws.Range('E30').Value = '  -3.64%  '
$ws.Range('B31').Value = 'Stellar'
$ws.Range('C31').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D31').Value = '0.08903'
$ws.Range('E31').Value = '  -1.57%  '
$ws.Range('B32').Value = 'ARBITRUM'
$ws.Range('C32').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D32').Value = '1.199'
$ws.Range('E32').Value = '  -1.45%  '
$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').Value = '0.7597'
$ws.Range('E33').Value = '  -1.38%  '
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').Value = '4.454'
$ws.Range('E34').Value = '  -2.28%  '
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').Value = '2.804'
$ws.Range('E35').Value = '  -7.54%  '
$ws.Range('B36').Value = 'Frax'
$ws.Range('C36').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D36').Value = '1.005'
$ws.Range('E36').Value = '  -1.32%  '
$ws.Range('B37').Value = 'TrustWalletToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D37').Value = '1.115'
$ws.Range('E37').Value = '  -2.00%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').Value = '0.01970'
$ws.Range('E38').Value = '  -0.64%  '
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').Value = '0.05279'
$ws.Range('E39').Value = '  -0.52%  '
$ws.Range('B40').Value = 'MXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D40').Value = '2.916'
$ws.Range('E40').Value = '  +1.62%  '
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').Value = '7.121'
$ws.Range('E41').Value = '  +2.47%  '
$ws.Range('B42').Value = 'Algorand'
$ws.Range('C42').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D42').Value = '0.1682'
$ws.Range('E42').Value = '  +0.32%  '
$ws.Range('B43').Value = 'TheSandbox'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D43').Value = '0.5032'
$ws.Range('E43').Value = '  -2.86%  '
$ws.Range('B44').Value = 'Aptos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D44').Value = '8.616'
$ws.Range('E44').Value = '  -1.20%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').Value = '10.57'
$ws.Range('E45').Value = '  -1.41%  '
$ws.Range('B46').Value = 'Quant'
$ws.Range('C46').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D46').Value = '106.73'
$ws.Range('E46').Value = '  -3.03%  '
$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D47').Value = '0.4719'
$ws.Range('E47').Value = '  -0.21%  '
$ws.Range('B48').Value = 'PaxDollar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D48').Value = '1.004'
$ws.Range('E48').Value = '  -1.40%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').Value = '0.06395'
$ws.Range('E49').Value = '  -1.14%  '
$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D50').Value = '1.666'
$ws.Range('E50').Value = '  -2.73%  '
$ws.Range('B51').Value = 'RenderToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D51').Value = '1.811'
$ws.Range('E51').Value = '  -2.20%  '
